$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two workers on rows 16 and 17 swap places (row 16 becomes DANIELA,
# row 17 becomes YAIR), and the "Salario Basico" (column G) figures are
# refreshed as part of the EC database update.

# Row 16 -> DANIELA MARTINEZ ALVAREZ
$ws.Range("C16").Value = "1007588337"
$ws.Range("D16").Value = "DANIELA MARTINEZ ALVAREZ"
$ws.Range("E16").Value = "2303"
$ws.Range("F16").Value = 7208
$ws.Range("G16").Value = 5406000

# Row 17 -> YAIR FERNEL OSORIO RAMIREZ
$ws.Range("C17").Value = "1143131446"
$ws.Range("D17").Value = "YAIR FERNEL OSORIO RAMIREZ"
$ws.Range("E17").Value = "2412"
$ws.Range("F17").Value = 388950
$ws.Range("G17").Value = 11802547
